# Update the predicted consumption values in column B (y_pred) while
# leaving the dates in column A and the headers in row 1 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    22896.25319999996,
    22564.97169999995,
    23076.73949999996,
    23002.79747499995,
    22295.87597499996,
    22233.39607499996
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
